# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Vega Modelo de Temuco - Plátano"
# as row 292, pushing the existing rows 292:320 down to 293:321.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 292, shifting rows 292:320 down to 293:321.
# -4121 = xlShiftDown
$ws.Rows("292:292").Insert(-4121)

# Populate the newly inserted row with the new weekly record.
$ws.Range("A292").Value = 10
$ws.Range("B292").Value = "Vega Modelo de Temuco"
$ws.Range("C292").Value = "La Araucanía"
$ws.Range("D292").Value = 44449
$ws.Range("E292").Value = 9
$ws.Range("F292").Value = "Fruta"
$ws.Range("G292").Value = 100108
$ws.Range("H292").Value = "Tropicales y subtropicales"
$ws.Range("I292").Value = 100108006
$ws.Range("J292").Value = "Plátano"
$ws.Range("K292").Value = "Sin especificar"
$ws.Range("L292").Value = "Pintón"
$ws.Range("M292").Value = 200
$ws.Range("N292").Value = 22000
$ws.Range("O292").Value = 22000
$ws.Range("P292").Value = 22000
$ws.Range("Q292").Value = "$/caja 20 kilos"
$ws.Range("R292").Value = "Ecuador"
$ws.Range("S292").Value = 1100
$ws.Range("T292").Value = 20

# Keep the date column formatted the same way as the rest of the column.
$ws.Range("D292").NumberFormat = "YYYY-MM-DD HH:MM:SS"
